$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title 1: "A Table, with a caption" -- merge each word with its trailing space
$title = $s.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$titleRange.Characters(1, 2).Text = "A "
$titleRange.Characters(3, 7).Text = "Table, "
$titleRange.Characters(10, 5).Text = "with "
$titleRange.Characters(15, 2).Text = "a "

# TextBox 3: "Demonstration of simple table syntax, with alignment"
$caption = $s.Shapes.Item(3)
$captionRange = $caption.TextFrame.TextRange
$captionRange.Characters(1, 14).Text = "Demonstration "
$captionRange.Characters(15, 3).Text = "of "
$captionRange.Characters(18, 7).Text = "simple "
$captionRange.Characters(25, 6).Text = "table "
$captionRange.Characters(31, 8).Text = "syntax, "
$captionRange.Characters(39, 5).Text = "with "
